$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-obsolete rows 8-10 (data was recomputed with new TPM; only 6 data rows remain)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Update remaining data rows (2-7) with recalculated TPM-based values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F2"
$ws.Range("C2").Value = "F2rl2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08786666666666666
$ws.Range("N2").Value = 0.2636
$ws.Range("O2").Value = 0.2208319231142997
$ws.Range("P2").Value = 0.2208319231142997
$ws.Range("Q2").Value = 0.003812798266666666
$ws.Range("R2").Value = 0.0343151844
$ws.Range("S2").Value = 0.01541454474042568
$ws.Range("T2").Value = 0.01541454474042568

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F2"
$ws.Range("C3").Value = "F2rl2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3100226666666667
$ws.Range("N3").Value = 0.9300680000000001
$ws.Range("O3").Value = 0.7791680768857003
$ws.Range("P3").Value = 0.7791680768857002
$ws.Range("Q3").Value = 0.01345281357466667
$ws.Range("R3").Value = 0.121075322172
$ws.Range("S3").Value = 0.05438761304111623
$ws.Range("T3").Value = 0.05438761304111622

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F2"
$ws.Range("C4").Value = "F2rl2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4205383333333333
$ws.Range("H4").Value = 1.261615
$ws.Range("I4").Value = 0.6764796878879081
$ws.Range("J4").Value = 0.6764796878879081
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08786666666666666
$ws.Range("N4").Value = 0.2636
$ws.Range("O4").Value = 0.2208319231142997
$ws.Range("P4").Value = 0.2208319231142997
$ws.Range("Q4").Value = 0.03695130155555555
$ws.Range("R4").Value = 0.332561714
$ws.Range("S4").Value = 0.149388310424048
$ws.Range("T4").Value = 0.149388310424048

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F2"
$ws.Range("C5").Value = "F2rl2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4205383333333333
$ws.Range("H5").Value = 1.261615
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3100226666666667
$ws.Range("N5").Value = 0.9300680000000001
$ws.Range("O5").Value = 0.7791680768857003
$ws.Range("P5").Value = 0.7791680768857002
$ws.Range("Q5").Value = 0.1303764155355556
$ws.Range("R5").Value = 1.17338773982
$ws.Range("S5").Value = 0.5270913774638601
$ws.Range("T5").Value = 0.5270913774638601

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F2"
$ws.Range("C6").Value = "F2rl2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1577256666666667
$ws.Range("H6").Value = 0.473177
$ws.Range("I6").Value = 0.2537181543305499
$ws.Range("J6").Value = 0.2537181543305499
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.08786666666666666
$ws.Range("N6").Value = 0.2636
$ws.Range("O6").Value = 0.2208319231142997
$ws.Range("P6").Value = 0.2208319231142997
$ws.Range("Q6").Value = 0.01385882857777778
$ws.Range("R6").Value = 0.1247294572
$ws.Range("S6").Value = 0.05602906794982605
$ws.Range("T6").Value = 0.05602906794982605

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F2"
$ws.Range("C7").Value = "F2rl2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1577256666666667
$ws.Range("H7").Value = 0.473177
$ws.Range("I7").Value = 0.2537181543305499
$ws.Range("J7").Value = 0.2537181543305499
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3100226666666667
$ws.Range("N7").Value = 0.9300680000000001
$ws.Range("O7").Value = 0.7791680768857003
$ws.Range("P7").Value = 0.7791680768857002
$ws.Range("Q7").Value = 0.04889853178177779
$ws.Range("R7").Value = 0.440086786036
$ws.Range("S7").Value = 0.1976890863807239
$ws.Range("T7").Value = 0.1976890863807239
